$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new resistor rows -----------------------------------------
# Before the edit, row 9 = "10KΩ" resistor and row 10 = "2K2Ω" resistor.
# After the edit, two new rows (220Ω, 360Ω) are inserted before the 10KΩ
# row, and one new row (220Ω, rich-text Ω) is inserted between the 10KΩ
# row and the 2K2Ω row.

# Insert two blank rows above the current row 9 (pushes "10KΩ" row to 11).
$ws.Rows(9).Insert()
$ws.Rows(9).Insert()

# Insert one blank row above the current row 12 (pushes "2K2Ω" row to 13).
$ws.Rows(12).Insert()

# --- Fill the three new rows ------------------------------------------------
# Row 9: Qtd 1, Encap (0805, Componente Resistor, Valor 220Ω (plain text)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "(0805"
$ws.Range("D9").Value = "Resistor"
$ws.Range("E9").Value = "220Ω"

# Row 10: Qtd 1, Encap (0805, Componente Resistor, Valor 360Ω (plain text)
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "(0805"
$ws.Range("D10").Value = "Resistor"
$ws.Range("E10").Value = "360Ω"

# Row 12: Qtd 3, Encap (0805, Componente Resistor, Valor 220Ω
# (this one is stored as rich text "220" + "Ω" runs in the original file,
# so split the formatting on the last character to reproduce that run split)
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "(0805"
$ws.Range("D12").Value = "Resistor"
$ws.Range("E12").Value = "220Ω"
$ws.Range("E12").Characters(4, 1).Font.Size = 11
$ws.Range("E12").Characters(4, 1).Font.Name = "Calibri"

# --- Update the selection shown in the saved view --------------------------
$ws.Range("H11").Select()
